$d = $word.ActiveDocument
$d.Content.Find.Execute("Devise and Bootstrap", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Devise-Bootstrap", 2)
